$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 for the new "Jurisdiction" metadata property.
# This shifts the existing rows 11-21 (Description ... Count) down to 12-22.
$ws.Rows.Item(11).Insert()

# The freshly inserted row doesn't carry the table's data-row formatting, so copy
# it over from the row directly below (which holds the shifted-down former row 11).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new "Jurisdiction" property row (value left blank, as in the source data).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh the "Date" metadata property with the new generation timestamp.
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"
